$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-03"

# Update the header label for the April column
$ws.Range("A5").Value = "April (through 04-03)"

# Update April row (row 5) values
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 7

# Update Total row (row 6) values
$ws.Range("B6").Value = 67
$ws.Range("C6").Value = 131
$ws.Range("D6").Value = 194
$ws.Range("E6").Value = 202
$ws.Range("F6").Value = 114
$ws.Range("G6").Value = 206
$ws.Range("H6").Value = 430
$ws.Range("I6").Value = 440
